$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# row 2 (context G=5489)
$ws.Range("H2").Value = 267.25
$ws.Range("I2").Value = 111.3
$ws.Range("J2").Value = 527.1667
$ws.Range("K2").Value = 111.3
$ws.Range("L2").Value = 527.1667
$ws.Range("M2").Value = 1.700000000000003
$ws.Range("N2").Value = -753.1667

# row 11 (context G=5533)
$ws.Range("H11").Value = 29.38889
$ws.Range("I11").Value = 29.38889
$ws.Range("K11").Value = 29.38889
$ws.Range("M11").Value = 110.61111

# row 40 (context G=5505)
$ws.Range("H40").Value = 5603.625
$ws.Range("J40").Value = 8599.200000000001
$ws.Range("L40").Value = 8599.200000000001
$ws.Range("N40").Value = -8949.200000000001

# row 98 (context G=36237)
$ws.Range("H98").Value = 1965.8334
$ws.Range("I98").Value = 1965.8334
$ws.Range("K98").Value = 1965.8334
$ws.Range("M98").Value = -467.8334

# row 111 (context G=27768)
$ws.Range("H111").Value = 1140.75
$ws.Range("I111").Value = 1017.1818
$ws.Range("K111").Value = 3051.5454
$ws.Range("M111").Value = 15.45460000000003

# row 122 (context G=36237)
$ws.Range("H122").Value = 1965.8334
$ws.Range("I122").Value = 1965.8334
$ws.Range("K122").Value = 5897.5002
$ws.Range("M122").Value = -3447.5002

# row 132 (context G=44049)
$ws.Range("H132").Value = 12562.435
$ws.Range("I132").Value = 11654.579
$ws.Range("J132").Value = 16874.75
$ws.Range("K132").Value = 34963.737
$ws.Range("L132").Value = 50624.25
$ws.Range("M132").Value = -32433.737
$ws.Range("N132").Value = -55684.25

# row 138 (context G=44169)
$ws.Range("H138").Value = 2173
$ws.Range("J138").Value = 2424.75
$ws.Range("L138").Value = 7274.25
$ws.Range("N138").Value = -17554.25

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# row 32 (context G=44147)
$ws.Range("H32").Value = 4479.769
$ws.Range("I32").Value = 4479.769
$ws.Range("K32").Value = 4479.769
$ws.Range("M32").Value = -4192.769

# row 86 (context G=10702)
$ws.Range("H86").Value = 0
$ws.Range("J86").Value = 0
$ws.Range("L86").Value = 0
$ws.Range("N86").Value = ""

# row 89 (context G=10702)
$ws.Range("H89").Value = 0
$ws.Range("J89").Value = 0
$ws.Range("L89").Value = 0
$ws.Range("N89").Value = ""

# row 110 (context G=27708)
$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = ""
$ws.Range("N110").Value = ""

# row 125 (context G=34251)
$ws.Range("H125").Value = 66888.89
$ws.Range("J125").Value = 66888.89
$ws.Range("L125").Value = 66888.89
$ws.Range("N125").Value = -76728.89

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# row 99 (context G=36198)
$ws.Range("H99").Value = 3119.5417
$ws.Range("I99").Value = 2812.2727
$ws.Range("K99").Value = 2812.2727
$ws.Range("M99").Value = -1314.2727

# row 122 (context G=36196)
$ws.Range("H122").Value = 1548.7587
$ws.Range("I122").Value = 1443.1578
$ws.Range("J122").Value = 1749.4
$ws.Range("K122").Value = 4329.4734
$ws.Range("L122").Value = 5248.200000000001
$ws.Range("M122").Value = -1879.4734
$ws.Range("N122").Value = -10148.2

# row 126 (context G=36198)
$ws.Range("H126").Value = 3119.5417
$ws.Range("I126").Value = 2812.2727
$ws.Range("K126").Value = 8436.8181
$ws.Range("M126").Value = -5966.8181

# row 134 (context G=44020)
$ws.Range("H134").Value = 1502.8611
$ws.Range("I134").Value = 1323.6
$ws.Range("K134").Value = 3970.8
$ws.Range("M134").Value = -1435.8

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# row 7 (context G=4728)
$ws.Range("H7").Value = 65.818184
$ws.Range("I7").Value = 74.666664
$ws.Range("J7").Value = 26
$ws.Range("K7").Value = 223.999992
$ws.Range("L7").Value = 78
$ws.Range("M7").Value = -111.999992
$ws.Range("N7").Value = -302

# row 92 (context G=19841)
$ws.Range("H92").Value = 2726.8333
$ws.Range("I92").Value = 1337
$ws.Range("J92").Value = 3190.111
$ws.Range("K92").Value = 4011
$ws.Range("L92").Value = 9570.332999999999
$ws.Range("M92").Value = -2763
$ws.Range("N92").Value = -12066.333

# row 109 (context G=27854)
$ws.Range("H109").Value = 502300
$ws.Range("I109").Value = 1000000
$ws.Range("K109").Value = 3000000
$ws.Range("M109").Value = -2998960

# row 118 (context G=27872)
$ws.Range("H118").Value = 249.6
$ws.Range("I118").Value = 249.6
$ws.Range("K118").Value = 748.8
$ws.Range("M118").Value = 494.2

# row 140 (context G=44097)
$ws.Range("H140").Value = 2000.9445
$ws.Range("I140").Value = 1307.8667
$ws.Range("J140").Value = 5466.3335
$ws.Range("K140").Value = 3923.6001
$ws.Range("L140").Value = 16399.0005
$ws.Range("M140").Value = 1256.3999
$ws.Range("N140").Value = -26759.0005

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# row 33 (context G=4450)
$ws.Range("H33").Value = 9244.75
$ws.Range("J33").Value = 9993
$ws.Range("L33").Value = 9993
$ws.Range("N33").Value = -10497

# row 96 (context G=18261)
$ws.Range("H96").Value = 52665
$ws.Range("J96").Value = 52665
$ws.Range("L96").Value = 52665
$ws.Range("N96").Value = -58157

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# row 40 (context G=36248)
$ws.Range("H40").Value = 5210.28
$ws.Range("I40").Value = 3728.842
$ws.Range("J40").Value = 9901.5
$ws.Range("K40").Value = 3728.842
$ws.Range("L40").Value = 9901.5
$ws.Range("M40").Value = -3592.842
$ws.Range("N40").Value = -10173.5

# row 105 (context G=18698)
$ws.Range("H105").Value = 46500
$ws.Range("J105").Value = 46500
$ws.Range("L105").Value = 46500
$ws.Range("N105").Value = -53488

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# row 33 (context G=2734)
$ws.Range("H33").Value = 11203.6
$ws.Range("I33").Value = 6509
$ws.Range("J33").Value = 14333.333
$ws.Range("K33").Value = 6509
$ws.Range("L33").Value = 14333.333
$ws.Range("M33").Value = -6259
$ws.Range("N33").Value = -14833.333

# row 36 (context G=2734)
$ws.Range("H36").Value = 11203.6
$ws.Range("I36").Value = 6509
$ws.Range("J36").Value = 14333.333
$ws.Range("K36").Value = 6509
$ws.Range("L36").Value = 14333.333
$ws.Range("M36").Value = -6259
$ws.Range("N36").Value = -14833.333

# row 40 (context G=3601)
$ws.Range("H40").Value = 67028
$ws.Range("J40").Value = 67028
$ws.Range("L40").Value = 67028
$ws.Range("N40").Value = -67326

# row 94 (context G=18075)
$ws.Range("H94").Value = 10594375
$ws.Range("J94").Value = 10594375
$ws.Range("L94").Value = 10594375
$ws.Range("N94").Value = -10596177

# row 126 (context G=36210)
$ws.Range("H126").Value = 3007.4482
$ws.Range("I126").Value = 1468.2632
$ws.Range("J126").Value = 5931.9
$ws.Range("K126").Value = 4404.7896
$ws.Range("L126").Value = 17795.7
$ws.Range("M126").Value = -1934.7896
$ws.Range("N126").Value = -22735.7

# row 132 (context G=44029)
$ws.Range("H132").Value = 1831.7561
$ws.Range("I132").Value = 1730.3158
$ws.Range("K132").Value = 5190.9474
$ws.Range("M132").Value = -2660.9474
